$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 1.389617666666667
$ws.Range("H2").Value = 4.168853
$ws.Range("I2").Value = 0.01005942246918146
$ws.Range("J2").Value = 0.01005942246918146
$ws.Range("M2").Value = 0.01989833333333333
$ws.Range("N2").Value = 0.059695
$ws.Range("O2").Value = 0.5455534129646046
$ws.Range("P2").Value = 0.5455534129646046
$ws.Range("Q2").Value = 0.02765107553722223
$ws.Range("R2").Value = 0.248859679835
$ws.Range("S2").Value = 0.005487952260514773
$ws.Range("T2").Value = 0.005487952260514773
$ws.Range("G3").Value = 1.389617666666667
$ws.Range("H3").Value = 4.168853
$ws.Range("I3").Value = 0.01005942246918146
$ws.Range("J3").Value = 0.01005942246918146
$ws.Range("O3").Value = 0.1997148627777118
$ws.Range("P3").Value = 0.1997148627777118
$ws.Range("Q3").Value = 0.01012243828988889
$ws.Range("R3").Value = 0.09110194460900001
$ws.Range("S3").Value = 0.002009016178055605
$ws.Range("T3").Value = 0.002009016178055605
$ws.Range("G4").Value = 1.389617666666667
$ws.Range("H4").Value = 4.168853
$ws.Range("I4").Value = 0.01005942246918146
$ws.Range("J4").Value = 0.01005942246918146
$ws.Range("M4").Value = 0.009290999999999999
$ws.Range("N4").Value = 0.027873
$ws.Range("O4").Value = 0.2547317242576836
$ws.Range("P4").Value = 0.2547317242576836
$ws.Range("Q4").Value = 0.012910937741
$ws.Range("R4").Value = 0.116198439669
$ws.Range("S4").Value = 0.002562454030611077
$ws.Range("T4").Value = 0.002562454030611078
$ws.Range("I5").Value = 0.01563117327130572
$ws.Range("J5").Value = 0.01563117327130572
$ws.Range("M5").Value = 0.01989833333333333
$ws.Range("N5").Value = 0.059695
$ws.Range("O5").Value = 0.5455534129646046
$ws.Range("P5").Value = 0.5455534129646046
$ws.Range("Q5").Value = 0.04296655739277778
$ws.Range("R5").Value = 0.386699016535
$ws.Range("S5").Value = 0.008527639926801936
$ws.Range("T5").Value = 0.008527639926801936
$ws.Range("I6").Value = 0.01563117327130572
$ws.Range("J6").Value = 0.01563117327130572
$ws.Range("O6").Value = 0.1997148627777118
$ws.Range("P6").Value = 0.1997148627777118
$ws.Range("S6").Value = 0.003121777624933457
$ws.Range("T6").Value = 0.003121777624933457
$ws.Range("I7").Value = 0.01563117327130572
$ws.Range("J7").Value = 0.01563117327130572
$ws.Range("M7").Value = 0.009290999999999999
$ws.Range("N7").Value = 0.027873
$ws.Range("O7").Value = 0.2547317242576836
$ws.Range("P7").Value = 0.2547317242576836
$ws.Range("Q7").Value = 0.020062096561
$ws.Range("R7").Value = 0.180558869049
$ws.Range("S7").Value = 0.003981755719570322
$ws.Range("T7").Value = 0.003981755719570323
$ws.Range("G8").Value = 0.6734466666666666
$ws.Range("H8").Value = 2.02034
$ws.Range("I8").Value = 0.004875070814774726
$ws.Range("J8").Value = 0.004875070814774726
$ws.Range("M8").Value = 0.01989833333333333
$ws.Range("N8").Value = 0.059695
$ws.Range("O8").Value = 0.5455534129646046
$ws.Range("P8").Value = 0.5455534129646046
$ws.Range("Q8").Value = 0.01340046625555556
$ws.Range("R8").Value = 0.1206041963
$ws.Range("S8").Value = 0.002659611521444487
$ws.Range("T8").Value = 0.002659611521444487
$ws.Range("G9").Value = 0.6734466666666666
$ws.Range("H9").Value = 2.02034
$ws.Range("I9").Value = 0.004875070814774726
$ws.Range("J9").Value = 0.004875070814774726
$ws.Range("O9").Value = 0.1997148627777118
$ws.Range("P9").Value = 0.1997148627777118
$ws.Range("Q9").Value = 0.004905610002222222
$ws.Range("R9").Value = 0.04415049002
$ws.Range("S9").Value = 0.0009736240988043619
$ws.Range("T9").Value = 0.000973624098804362
$ws.Range("G10").Value = 0.6734466666666666
$ws.Range("H10").Value = 2.02034
$ws.Range("I10").Value = 0.004875070814774726
$ws.Range("J10").Value = 0.004875070814774726
$ws.Range("M10").Value = 0.009290999999999999
$ws.Range("N10").Value = 0.027873
$ws.Range("O10").Value = 0.2547317242576836
$ws.Range("P10").Value = 0.2547317242576836
$ws.Range("Q10").Value = 0.006256992979999999
$ws.Range("R10").Value = 0.05631293682
$ws.Range("S10").Value = 0.001241835194525876
$ws.Range("T10").Value = 0.001241835194525876
$ws.Range("G11").Value = 133.9185306666667
$ws.Range("H11").Value = 401.755592
$ws.Range("I11").Value = 0.9694343334447382
$ws.Range("J11").Value = 0.9694343334447382
$ws.Range("M11").Value = 0.01989833333333333
$ws.Range("N11").Value = 0.059695
$ws.Range("O11").Value = 0.5455534129646046
$ws.Range("P11").Value = 0.5455534129646046
$ws.Range("Q11").Value = 2.664755562715555
$ws.Range("R11").Value = 23.98280006444
$ws.Range("S11").Value = 0.5288782092558434
$ws.Range("T11").Value = 0.5288782092558434
$ws.Range("G12").Value = 133.9185306666667
$ws.Range("H12").Value = 401.755592
$ws.Range("I12").Value = 0.9694343334447382
$ws.Range("J12").Value = 0.9694343334447382
$ws.Range("O12").Value = 0.1997148627777118
$ws.Range("P12").Value = 0.1997148627777118
$ws.Range("Q12").Value = 0.9755072168862223
$ws.Range("R12").Value = 8.779564951975999
$ws.Range("S12").Value = 0.1936104448759184
$ws.Range("T12").Value = 0.1936104448759184
$ws.Range("G13").Value = 133.9185306666667
$ws.Range("H13").Value = 401.755592
$ws.Range("I13").Value = 0.9694343334447382
$ws.Range("J13").Value = 0.9694343334447382
$ws.Range("M13").Value = 0.009290999999999999
$ws.Range("N13").Value = 0.027873
$ws.Range("O13").Value = 0.2547317242576836
$ws.Range("P13").Value = 0.2547317242576836
$ws.Range("Q13").Value = 1.244237068424
$ws.Range("R13").Value = 11.198133615816
$ws.Range("S13").Value = 0.2469456793129763
$ws.Range("T13").Value = 0.2469456793129764
